$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 13 with the new Narcotics Stock report test case.
# New shared strings must be appended in this order: TC014, then the path.
$ws.Cells.Item(13, 5).Value = "TC014"
$ws.Cells.Item(13, 1).Value = "Pharmacy\Reports\Stock\TC04NarcoticStockReport.py"
$ws.Cells.Item(13, 2).Value = "Norun"
$ws.Cells.Item(13, 3).Value = "PharmacyReport"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 8).Value = "SNCH"

# Copy styles from row 12 to row 13 so formatting (borders, fill, alignment)
# matches the rest of the table, including the special style on column A.
$ws.Range("A12:H12").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the target state
$ws.Range("B17").Select()
